# Apply crypto price/volume updates (and a 2-row coin swap) as produced by the
# "Updated cryptos list" GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.783.51"
$ws.Range("E2").Value = "  -7.75%  "
$ws.Range("D3").Value = "2.520.57"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.56"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.15"
$ws.Range("E6").Value = "  -7.05%  "
$ws.Range("E7").Value = "  -5.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -6.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("E10").Value = "  -8.17%  "
$ws.Range("E11").Value = "  -5.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("E12").Value = "  -6.24%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "2.906.30"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "2.520.91"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("E16").Value = "  -6.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  -6.56%  "
$ws.Range("D18").Value = "42.793.31"
$ws.Range("E18").Value = "  -8.18%  "
$ws.Range("D19").Value = "0.0₃0961"
$ws.Range("E19").Value = "  -5.16%  "
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("E21").Value = "  -5.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.50"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.60"
$ws.Range("E23").Value = "  -6.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.15"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.96"
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.73"
$ws.Range("E30").Value = "  -5.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.47"
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.25"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.75"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0799"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.91"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.38"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0307"
$ws.Range("E42").Value = "  -7.34%  "
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("D44").Value = "2.016.36"
$ws.Range("E44").Value = "  -5.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.64"
$ws.Range("E46").Value = "  -9.26%  "
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.85"
$ws.Range("E48").Value = "  -7.24%  "
$ws.Range("D49").Value = "2.762.87"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.68"
$ws.Range("E50").Value = "  -6.55%  "
$ws.Range("E51").Value = "  -8.50%  "

Write-Host "Update complete"
